# Price update for 2026-02-07
# Append a new tracking row (row 38) to "Sheet 1" with the latest scrape:
#   A38 Date       = 2026-02-07
#   B38 Price      = 138000
#   C38 Discount   = 0
#   D38 Incredible = 0
#
# Every column in this sheet stores its values as plain text (shared
# strings) even though Price/Discount/Incredible look numeric and Date
# looks like a date. A bare "$range.Value = ..." would let Excel's
# smart-typing reinterpret "2026-02-07" as a date serial and
# "138000"/"0" as numbers, so the new cells are switched to Text format
# before the values are typed in, then restored to the workbook's
# default (Normal) style once the literal text is safely stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("A38:D38")

# Force Text number format so the literal strings aren't reinterpreted
# as a date / numbers by Excel's input parser.
$target.NumberFormat = "@"

$ws.Range("A38").Value = "2026-02-07"
$ws.Range("B38").Value = "138000"
$ws.Range("C38").Value = "0"
$ws.Range("D38").Value = "0"

# Restore the default (Normal) cell style now that the text is stored,
# matching the rest of the sheet which carries no explicit style.
$target.Style = "Normal"
